$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Accredited)
$ws.Range("C2").Value = 0.982
$ws.Range("D2").Value = 6.416
$ws.Range("E2").Value = 2.538
$ws.Range("F2").Value = 2.514
$ws.Range("G2").Value = 3.807
$ws.Range("H2").Value = 3.136
$ws.Range("I2").Value = 1.508
$ws.Range("J2").Value = 3.328
$ws.Range("K2").Value = 5.65
$ws.Range("L2").Value = 4.884
$ws.Range("M2").Value = 3.017
$ws.Range("N2").Value = 3.352
$ws.Range("O2").Value = 41.132

# Row 3 (Unaccredited)
$ws.Range("C3").Value = 1.58
$ws.Range("D3").Value = 8.475
$ws.Range("E3").Value = 3.663
$ws.Range("F3").Value = 4.094
$ws.Range("G3").Value = 6.105
$ws.Range("H3").Value = 4.453
$ws.Range("I3").Value = 1.748
$ws.Range("J3").Value = 4.453
$ws.Range("K3").Value = 8.451000000000001
$ws.Range("L3").Value = 8.571
$ws.Range("M3").Value = 3.759
$ws.Range("N3").Value = 3.519
$ws.Range("O3").Value = 58.871

# Row 4 (COL_TOT)
$ws.Range("C4").Value = 2.562
$ws.Range("D4").Value = 14.891
$ws.Range("E4").Value = 6.201
$ws.Range("F4").Value = 6.608000000000001
$ws.Range("G4").Value = 9.912000000000001
$ws.Range("H4").Value = 7.589
$ws.Range("I4").Value = 3.256
$ws.Range("J4").Value = 7.781000000000001
$ws.Range("K4").Value = 14.101
$ws.Range("L4").Value = 13.455
$ws.Range("M4").Value = 6.776
$ws.Range("N4").Value = 6.871
$ws.Range("O4").Value = 100.003
